$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Status" column header and values
$ws.Range("E1").Value = "Status"
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 1

# Update selection to match the new last cell used (as captured in the source workbook)
$ws.Range("E4").Select()
